# "Generate Report for Handback"
# Update the localization-status report: the aaaa355a-... file has now been
# handed back (in sync with en-US) instead of merely "Ready for handoff",
# so refresh its status + handback timestamps on every sheet, and clear the
# stale "handback file not latest" error now that it is resolved.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns for the aaaa355a row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: Status, Latest Handback DateTime, Error Detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-26 15:00:52"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8

# --- de-de sheet: Status, Latest Handback DateTime, Error Detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-26 15:01:19"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
